# Auto-generated edit script: update cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "51.462.07"
$ws.Range("E2").Value = "  +0.66%  "

# Row 3
$ws.Range("D3").Value = "2.982.49"
$ws.Range("E3").Value = "  +1.24%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.40"
$ws.Range("E5").Value = "  +1.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.77"
$ws.Range("E6").Value = "  +2.27%  "

# Row 7
$ws.Range("E7").Value = "  +0.97%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.66"
$ws.Range("E10").Value = "  +0.54%  "

# Row 11
$ws.Range("E11").Value = "  -0.94%  "

# Row 12
$ws.Range("E12").Value = "  +0.46%  "

# Row 13
$ws.Range("D13").Value = "3.449.87"
$ws.Range("E13").Value = "  +1.23%  "

# Row 14
$ws.Range("E14").Value = "  +1.59%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.79"
$ws.Range("E15").Value = "  +2.26%  "

# Row 16
$ws.Range("D16").Value = "2.975.14"
$ws.Range("E16").Value = "  +0.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.20"
$ws.Range("E17").Value = "  +1.59%  "

# Row 18
$ws.Range("E18").Value = "  -0.10%  "

# Row 19
$ws.Range("D19").Value = "51.457.19"
$ws.Range("E19").Value = "  +0.69%  "

# Row 20
$ws.Range("E20").Value = "  +0.24%  "

# Row 21
$ws.Range("E21").Value = "  +0.76%  "

# Row 22
$ws.Range("E22").Value = "  +0.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.31"
$ws.Range("E23").Value = "  +2.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.81"
$ws.Range("E24").Value = "  +0.38%  "

# Row 25
$ws.Range("E25").Value = "  +2.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.83"
$ws.Range("E26").Value = "  -4.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.34"
$ws.Range("E27").Value = "  -1.71%  "

# Row 28
$ws.Range("E28").Value = "  +2.36%  "

# Row 29
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.10"
$ws.Range("E30").Value = "  +1.68%  "

# Row 31
$ws.Range("E31").Value = "  -0.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.39"
$ws.Range("E32").Value = "  +3.69%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.75"
$ws.Range("E33").Value = "  +3.69%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.44"
$ws.Range("E34").Value = "  +1.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  +0.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0442"
$ws.Range("E36").Value = "  -0.37%  "

# Row 38
$ws.Range("E38").Value = "  +3.75%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.81"
$ws.Range("E39").Value = "  +2.13%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.57"
$ws.Range("E40").Value = "  +3.30%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  +0.53%  "

# Row 42
$ws.Range("E42").Value = "  +2.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.21"
$ws.Range("E43").Value = "  +5.63%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.79"
$ws.Range("E44").Value = "  +11.89%  "

# Row 45
$ws.Range("E45").Value = "  +0.49%  "

# Row 46
$ws.Range("E46").Value = "  -0.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.273"
$ws.Range("E47").Value = "  -0.02%  "

# Row 48
$ws.Range("E48").Value = "  +2.40%  "

# Row 49
$ws.Range("D49").Value = "2.027.65"
$ws.Range("E49").Value = "  +1.84%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "3.279.34"
$ws.Range("E50").Value = "  +1.24%  "

# Row 51
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0334"
$ws.Range("E51").Value = "  +1.85%  "
